# Cambios del correo del 16 de feb
#
# Adds two new report rows (9 and 10) mirroring row 8's structure, fills in
# the "Tipo de documento" / "Denominación" / hyperlink columns for rows 8-10,
# clears the long "Nota" text from column K, widens two of the new rows,
# narrows column K, extends the data-validation range, and leaves the
# "Hidden_1" helper sheet as the active tab (mirroring the author's final
# on-screen state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- duplicate row 8 into rows 9 and 10, preserving its formatting -------
$ws.Rows.Item(8).Copy()
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(8).Copy()
$ws.Rows.Item(9).Insert()

# Row heights for the three data rows (content got shorter, so the tall
# 195pt row shrinks down to 45pt for all three entries).
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 45

$nl = [char]10

# --- row 8: Contable / Estado de Situación Financiera --------------------
$ws.Range("D8").Value = "Contable"
$ws.Range("E8").Value = "Estado de Situación Financiera"
$ws.Range("H8").Value = "Secretaría Admnistrativa (UPP)"
$ws.Range("K8").ClearContents()

$ws.Hyperlinks.Add(
    $ws.Range("F8"),
    "http://www.upp.edu.mx/leygralcontabilidad/mc/01-edosfin/2020/a_diciembre_2020/estado-de-situacion-financiera.pdf",
    "",
    "",
    "http://www.upp.edu.mx/leygralcontabilidad/mc/01-edosfin/2020/a_diciembre_2020/estado-de-situacion-financiera.pdf"
)
$ws.Hyperlinks.Add(
    $ws.Range("G8"),
    "http://transparencia.hidalgo.gob.mx/descargables/dependencias/finanzasadmon/16edofinanciero.pdf",
    "",
    "",
    "http://transparencia.hidalgo.gob.mx/" + $nl + "descargables/dependencias/finanza" + $nl + "sadmon/16edofinanciero.pdf"
)

# --- row 9: Presupuestal / Presupuesto Anual de Egresos Modificado -------
$ws.Range("D9").Value = "Presupuestal"
$ws.Range("E9").Value = "Presupuesto Anual de Egresos Modificado"
$ws.Range("H9").Value = "Subdirección de Programación y Presupuesto (UPP)"
$ws.Range("K9").ClearContents()

$ws.Hyperlinks.Add(
    $ws.Range("F9"),
    "http://www.upp.edu.mx/leygralcontabilidad/mc/02-edospres/02-programa-anual-de-egresos-autorizados/2020-Presupuesto_Anual_de_Egresos_Resumen.pdf",
    "",
    "",
    "http://www.upp.edu.mx/leygralcontabilidad/mc/02-edospres/02-programa-anual-de-egresos-autorizados/2020-Presupuesto_Anual_de_Egresos_Resumen.pdf"
)
$ws.Hyperlinks.Add(
    $ws.Range("G9"),
    "http://transparencia.hidalgo.gob.mx/descargables/dependencias/finanzasadmon/16edofinanciero.pdf",
    "",
    "",
    "http://transparencia.hidalgo.gob.mx/" + $nl + "descargables/dependencias/finanza" + $nl + "sadmon/16edofinanciero.pdf"
)

# --- row 10: Programático / Adecuaciones Prespuestarias POA --------------
$ws.Range("D10").Value = "Programático"
$ws.Range("E10").Value = "Adecuaciones Prespuestarias Programa Operativo Anual"
$ws.Range("H10").Value = "Subdirección de Programación y Presupuesto (UPP)"
$ws.Range("K10").ClearContents()

$ws.Hyperlinks.Add(
    $ws.Range("F10"),
    "http://www.upp.edu.mx/leygralcontabilidad/mc/02-edospres/03-programa-operativo-anual-autorizado/2020-POA_Programatico.pdf",
    "",
    "",
    "http://www.upp.edu.mx/leygralcontabilidad/mc/02-edospres/03-programa-operativo-anual-autorizado/2020-POA_Programatico.pdf"
)
$ws.Hyperlinks.Add(
    $ws.Range("G10"),
    "http://transparencia.hidalgo.gob.mx/descargables/dependencias/finanzasadmon/16edofinanciero.pdf",
    "",
    "",
    "http://transparencia.hidalgo.gob.mx/" + $nl + "descargables/dependencias/finanza" + $nl + "sadmon/16edofinanciero.pdf"
)

# --- column width: K shrinks now that the long note text is gone ---------
$ws.Columns.Item(11).ColumnWidth = 21.14

# --- extend the "Tipo de documento" dropdown down through row 201 --------
$ws.Range("D8:D201").Validation.Delete()
$ws.Range("D8:D201").Validation.Add(3, 1, 1, "=Hidden_13")
$ws.Range("D8:D201").Validation.IgnoreBlank = $true
$ws.Range("D8:D201").Validation.InCellDropdown = $true
$ws.Range("D8:D201").Validation.ShowInput = $false
$ws.Range("D8:D201").Validation.ShowError = $true

# --- view/selection state on the main sheet -------------------------------
$ws.Activate()
$ws.Range("A8:J10").Select()

# --- Hidden_1 ends up as the active/visible tab ---------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
